$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.012.10'
$ws.Range("E2").Value = '  +2.71%  '
$ws.Range("D3").Value = '3.087.75'
$ws.Range("E3").Value = '  +4.64%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").Value = "'580.17"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.07%  '
$ws.Range("D6").Value = "'168.62"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +5.81%  '
$ws.Range("D7").Value = "'1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.06%  '
$ws.Range("D8").Value = '3.082.61'
$ws.Range("E8").Value = '  +4.61%  '
$ws.Range("D9").Value = "'0.524"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.00%  '
$ws.Range("E10").Value = '  -0.56%  '
$ws.Range("E11").Value = '  +2.38%  '
$ws.Range("D12").Value = "'0.483"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +5.16%  '
$ws.Range("E13").Value = '  +1.84%  '
$ws.Range("D14").Value = "'36.56"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +6.44%  '
$ws.Range("E15").Value = '  -0.52%  '
$ws.Range("D16").Value = '3.598.24'
$ws.Range("E16").Value = '  +4.57%  '
$ws.Range("D17").Value = '66.966.15'
$ws.Range("E17").Value = '  +2.53%  '
$ws.Range("D18").Value = "'7.21"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +3.18%  '
$ws.Range("D19").Value = '3.084.49'
$ws.Range("E19").Value = '  +4.22%  '
$ws.Range("D20").Value = "'16.14"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +7.82%  '
$ws.Range("D21").Value = "'465.81"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +4.55%  '
$ws.Range("D22").Value = "'0.717"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +4.24%  '
$ws.Range("D23").Value = "'7.54"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +4.07%  '
$ws.Range("D24").Value = "'83.30"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.30%  '
$ws.Range("E25").Value = '  +6.52%  '
$ws.Range("D26").Value = "'12.89"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +6.57%  '
$ws.Range("E27").Value = '  +1.35%  '
$ws.Range("E29").Value = '  -0.36%  '
$ws.Range("E30").Value = '  +0.51%  '
$ws.Range("E31").Value = '  +3.43%  '
$ws.Range("E32").Value = '  +0.90%  '
$ws.Range("E33").Value = '  +3.45%  '
$ws.Range("E34").Value = '  +3.22%  '
$ws.Range("E35").Value = '  -0.06%  '
$ws.Range("E36").Value = '  +2.80%  '
$ws.Range("D37").Value = "'5.88"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.08%  '
$ws.Range("D38").Value = "'2.13"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +7.17%  '
$ws.Range("E39").Value = '  +5.64%  '
$ws.Range("D40").Value = "'0.319"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +7.00%  '
$ws.Range("D41").Value = "'50.21"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.42%  '
$ws.Range("E42").Value = '  +1.91%  '
$ws.Range("E43").Value = '  +2.46%  '
$ws.Range("E44").Value = '  -0.76%  '
$ws.Range("D45").Value = "'0.0360"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.46%  '
$ws.Range("D46").Value = "'384.04"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.28%  '
$ws.Range("D47").Value = '2.764.77'
$ws.Range("E47").Value = '  +2.00%  '
$ws.Range("D48").Value = "'134.49"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.13%  '
$ws.Range("E49").Value = '  +0.02%  '
$ws.Range("D50").Value = "'24.65"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +6.40%  '
$ws.Range("E51").Value = '  +2.56%  '
